$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13467.6364916629
$ws.Range("C2").Value = 13088.8344317158
$ws.Range("E2").Value = 9272.5924893091
$ws.Range("F2").Value = -1.43471162396448

$ws.Range("B3").Value = 13050.666638977
$ws.Range("C3").Value = 12456.3679112785
$ws.Range("E3").Value = 8926.87450595342
$ws.Range("F3").Value = 387.807600717996

$ws.Range("B4").Value = 12578.722509054
$ws.Range("C4").Value = 11935.7096041208
$ws.Range("E4").Value = 8556.76838299946
$ws.Range("F4").Value = 350.692416130012

$ws.Range("B5").Value = 12274.9658789007
$ws.Range("C5").Value = 11680.3148669912
$ws.Range("E5").Value = 8331.01356064185
$ws.Range("F5").Value = 330.644517818044

$ws.Range("B6").Value = 11932.2251655432
$ws.Range("C6").Value = 10676.4209242692
$ws.Range("E6").Value = 8082.03939317868
$ws.Range("F6").Value = 278.441679893662

$ws.Range("B7").Value = 5250.94347634098
$ws.Range("C7").Value = 7531.82370497339
$ws.Range("E7").Value = 8472.0050600299
$ws.Range("F7").Value = 163.665365208471
